$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 5790
$ws1.Range("F8").Value = 423
$ws1.Range("F9").Value = 3923
$ws1.Range("F16").Value = 113
$ws1.Range("F18").Value = 638
$ws1.Range("F19").Value = 3934
$ws1.Range("F20").Value = 141
$ws1.Range("F22").Value = 5428
$ws1.Range("F24").Value = 2143
$ws1.Range("F26").Value = 369
$ws1.Range("F27").Value = 8090
$ws1.Range("F29").Value = 11
$ws1.Range("G29").Value = 44.1
$ws1.Range("F30").Value = 2215
$ws1.Range("F31").Value = 2229
$ws1.Range("F32").Value = 1343
$ws1.Range("F33").Value = 176
$ws1.Range("F34").Value = 1329
$ws1.Range("F46").Value = 2145
$ws1.Range("F47").Value = 140
$ws1.Range("F48").Value = 236
$ws1.Range("F49").Value = 1224

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 130

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 593
$ws3.Range("F3").Value = 779

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 593
$ws4.Range("F6").Value = 779
$ws4.Range("F7").Value = 5790
$ws4.Range("F8").Value = 423
$ws4.Range("F9").Value = 3923
$ws4.Range("F15").Value = 113
$ws4.Range("F18").Value = 638
$ws4.Range("F19").Value = 3934
$ws4.Range("F21").Value = 141
$ws4.Range("F23").Value = 5428
$ws4.Range("F25").Value = 2143
$ws4.Range("F27").Value = 369
$ws4.Range("F28").Value = 8090
$ws4.Range("F30").Value = 2215
$ws4.Range("F31").Value = 2229
$ws4.Range("F32").Value = 176
$ws4.Range("F33").Value = 1329
$ws4.Range("F44").Value = 2146
$ws4.Range("F45").Value = 140
$ws4.Range("F47").Value = 236
$ws4.Range("F49").Value = 1224
